$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new title row above the existing table, shifting rows 1-13
#    down to rows 2-14.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).Insert()

# ---------------------------------------------------------------------------
# 2. Fix up / clarify several label strings in column A (rows now 3,6,8,9,
#    10,11 after the shift).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value2  = "Мощность двигателя (лошадиные силы):"
$ws.Range("A6").Value2  = "Период страховки (для иностранных агентов):"
$ws.Range("A8").Value2  = "Возраст водителя (период в годовом выражение):"
$ws.Range("A9").Value2  = "Водительский стаж (количество полных лет):"
$ws.Range("A10").Value2 = "ТС зарегестрировано в иностраном государстве:"
$ws.Range("A11").Value2 = "Юридическая форма:"

# ---------------------------------------------------------------------------
# 3. New title row text (A1) - B1 stays empty.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "Результаты расчета"

# ---------------------------------------------------------------------------
# 4. Column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 59.1
$ws.Columns.Item(2).ColumnWidth = 59.1

# ---------------------------------------------------------------------------
# 5. Body formatting (rows 2-13): left/bottom aligned text, thin borders -
#    column A gets left+right+bottom, column B gets right+bottom only (the
#    shared vertical divider comes from column A's right edge).
# ---------------------------------------------------------------------------
$bodyA = $ws.Range("A2:A13")
$bodyA.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$bodyA.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignBottom
$bodyA.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
$bodyA.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
$bodyA.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

$bodyB = $ws.Range("B2:B13")
$bodyB.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$bodyB.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignBottom
$bodyB.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
$bodyB.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

# ---------------------------------------------------------------------------
# 6. Title row (A1:B1) formatting: blue font, thick outer box.
# ---------------------------------------------------------------------------
$title = $ws.Range("A1:B1")
$title.Font.Color = 16711680
$title.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$title.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignBottom

$titleA = $ws.Range("A1")
$titleA.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThick
$titleA.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThick
$titleA.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThick

$titleB = $ws.Range("B1")
$titleB.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThick
$titleB.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThick
$titleB.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThick

# ---------------------------------------------------------------------------
# 7. Footer row (A14:B14) formatting: same blue font + thick box as title.
# ---------------------------------------------------------------------------
$footer = $ws.Range("A14:B14")
$footer.Font.Color = 16711680
$footer.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$footer.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignBottom

$footerA = $ws.Range("A14")
$footerA.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThick
$footerA.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThick
$footerA.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThick

$footerB = $ws.Range("B14")
$footerB.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThick
$footerB.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThick
$footerB.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThick

# ---------------------------------------------------------------------------
# 8. B9 ("Водительский стаж ...") value is a plain number, not text.
# ---------------------------------------------------------------------------
$ws.Range("B9").Value2 = 0

# ---------------------------------------------------------------------------
# 9. Gridlines stay visible; selection highlights the new title row.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("A1:B1").Select()
